$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# New rows of Issue / Notes data appended to the table
$newRows = @(
    @("Nominal Effort could be improved as an index of exploitation rate", "Can we derive effort / habitat area. There is the potential to borrow information on catchability among areas/models - priors, metaanalysis, EM. ", 28.8),
    @("Catches are expanded to totals using expansion factor - no uncertainty", "How can we get observation error in total catches? How are expansion factors calculated - can we do bootstrapping etc?", 28.8),
    @("Discard mortality rate assumed to be 9% but from a study elsehwere", "Lyle et al. 2006.  This is used to include discard mortality in total catch data (in model conditioning [Catch = ExpFac x (Kept + Rel * DiscMort)] and used in projections that would affect any kind of regulation affecting discarding such as size limits, bag limits etc. ", 57.6),
    @("Total recreational effort", "currently calculated by Duration_hrs x Npersons x ExpWt (what is the 'expansion factor'??)", 28.8)
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2)).RowHeight = $newRows[$i][2]
}

# Selection now sits on the newly noted "total recreational effort" cell
$ws.Range("B17:B18").Select()

# Resize/reposition the workbook window to match the author's session
$win = $excel.ActiveWindow
$win.Left = 3234
$win.Top = 936
$win.Width = 13422
$win.Height = 12246

$wb.Application.CutCopyMode = $false
